# Update for next seminar:
#   - bump the cached "datetimeFigureOut" date field (slide master + every
#     slide layout) from 2020-08-24 to 2020-08-31
#   - swap the moderator name "Kyle" -> "David" on the rules/instructions
#     slide (only the standalone, red-highlighted "Kyle" run - not the
#     "Kyle Murphy" credit further down the deck)

$p = $ppt.ActivePresentation

$oldDate = "2020-08-24"
$newDate = "2020-08-31"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePh = $false
        if ($shp.Type -eq 14) {
            try {
                if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePh = $true
                }
            } catch {
                $isDatePh = $false
            }
        }
        if ($isDatePh -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every custom (slide) layout has its own date placeholder too.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Rules slide: "Questions should be asked to Kyle." -> "...David."
$slide = $p.Slides.Item(1)
$shape = $slide.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$fullText = $tr.Text
$idx = $fullText.IndexOf("Kyle")
if ($idx -ge 0) {
    $chars = $tr.Characters($idx + 1, 4)
    $chars.Text = "David"
}
